$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 9: downsample percent (C9) 5 -> 10, and add a note in J9 ---
$ws.Cells.Item(9, 3).Value = 10
$ws.Range("J9").Value = "this one showed more values beyond lenthscale5 to train on"

# --- Row 10 (hybrid_40_lml.csv): fill in the hyperparameter columns ---
$ws.Cells.Item(10, 2).Value = 0.5
$ws.Cells.Item(10, 3).Value = 10
$ws.Cells.Item(10, 4).Value = 50
$ws.Cells.Item(10, 5).Value = 0.01
$ws.Cells.Item(10, 6).Value = 0.5
$ws.Cells.Item(10, 7).Value = 50
$ws.Cells.Item(10, 8).Value = 0.004

# --- Row 11 (hybrid_30_lml.csv): fill in the hyperparameter columns ---
$ws.Cells.Item(11, 2).Value = 0.5
$ws.Cells.Item(11, 3).Value = 10
$ws.Cells.Item(11, 4).Value = 50
$ws.Cells.Item(11, 5).Value = 0.01
$ws.Cells.Item(11, 6).Value = 0.5
$ws.Cells.Item(11, 7).Value = 50
$ws.Cells.Item(11, 8).Value = 0.004

# --- Row 12 (hybrid_20_lml.csv): fill in the hyperparameter columns ---
$ws.Cells.Item(12, 2).Value = 0.5
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 50
$ws.Cells.Item(12, 5).Value = 0.01
$ws.Cells.Item(12, 6).Value = 0.5
$ws.Cells.Item(12, 7).Value = 50
$ws.Cells.Item(12, 8).Value = 0.004

# --- Insert a brand new row 13 for the new "hybrid_10_lml_big.csv" entry ---
# (this pushes the former row 13, hybrid_10_lml.csv, down to row 14)
$ws.Rows.Item(13).Insert()

# New row 13 content
$ws.Range("A13").Value = "hybrid_10_lml_big.csv"
$ws.Cells.Item(13, 2).Value = 0.5
$ws.Cells.Item(13, 3).Value = 20
$ws.Cells.Item(13, 4).Value = 50
$ws.Cells.Item(13, 5).Value = 0.01
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 50
$ws.Cells.Item(13, 8).Value = 0.004

# Copy the bottom-border formatting (style used by A4/A14) onto A13 so it
# matches the rest of the "last row in block" styling without minting a new
# style entry.
$ws.Range("A4").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# --- Row 14 (the original hybrid_10_lml.csv row, shifted down by the insert) ---
# Fill in its hyperparameter values and drop the carried-over border
# formatting on B:H to match the plain (unstyled) value cells.
$ws.Range("B14:H14").ClearFormats()
$ws.Cells.Item(14, 2).Value = 0.5
$ws.Cells.Item(14, 3).Value = 10
$ws.Cells.Item(14, 4).Value = 50
$ws.Cells.Item(14, 5).Value = 0.01
$ws.Cells.Item(14, 6).Value = 0.5
$ws.Cells.Item(14, 7).Value = 50
$ws.Cells.Item(14, 8).Value = 0.004

# --- Update the selection to reflect where the user ended up (A15) ---
$ws.Range("A15").Select()
